$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in rows 15-17 that actually carry data which needs to rotate.
# (Columns such as D, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AY, I are identical
# across the three rows, so they are intentionally left untouched.)
$cols = @("A","B","E","F","G","H","P","Q","R","AC","AW","AX")

# Capture current ("before") values of rows 15, 16 and 17 before overwriting anything,
# because the update is a rotation: new15 = old16, new16 = old17, new17 = old15.
$row15 = @{}
$row16 = @{}
$row17 = @{}
foreach ($col in $cols) {
    $row15[$col] = $ws.Range("${col}15").Value2
    $row16[$col] = $ws.Range("${col}16").Value2
    $row17[$col] = $ws.Range("${col}17").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}15").Value2 = $row16[$col]
    $ws.Range("${col}16").Value2 = $row17[$col]
    $ws.Range("${col}17").Value2 = $row15[$col]
}
